$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates per diff (only row 2 is affected)
$ws.Range("C2").Value = 0.1741947494915157
$ws.Range("E2").Value = 0.09111564259528658
$ws.Range("G2").Value = 0.07877662572384592

$ws.Range("H2").Value = "태양계"
$ws.Range("I2").Value = 0.2371271638948406
$ws.Range("J2").Value = "행성"
$ws.Range("K2").Value = 0.2188866128260067
$ws.Range("L2").Value = "위성"
$ws.Range("M2").Value = 0.2096034752284752

$ws.Range("O2").Value = 0.2053872053872054
$ws.Range("Q2").Value = 0.1750841750841751
$ws.Range("R2").Value = "목성"
$ws.Range("S2").Value = 0.1414141414141414
